# Normalise the "Recorded By" column (G): whenever the comma-separated
# list of recorders currently starts with the literal entry "System",
# flip the whole list end-to-end so "System" lands last instead of first.
# Rows whose list does not begin with "System" (e.g. already-reordered
# rows, or rows with only a single entry) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ', '

        if ($parts.Count -gt 1 -and $parts[0] -eq "System") {
            $n = $parts.Count
            $reversedParts = @()
            for ($i = $n - 1; $i -ge 0; $i--) {
                $reversedParts += $parts[$i]
            }
            $newVal = [string]::Join(", ", $reversedParts)
            $cell.Value = $newVal
        }
    }
}
